$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 112044176
$ws.Range("Q9").Value = 554725
$ws.Range("R9").Value = 6697571
$ws.Range("Z9").ClearContents() | Out-Null
$ws.Range("AB9").ClearContents() | Out-Null

# Row 10
$ws.Range("Q10").Value = 554722
$ws.Range("R10").Value = 6697604
$ws.Range("Z10").ClearContents() | Out-Null
$ws.Range("AB10").ClearContents() | Out-Null

# Row 11
$ws.Range("A11").Value = 112044178
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("Q11").Value = 554737
$ws.Range("R11").Value = 6697621
$ws.Range("Z11").ClearContents() | Out-Null
$ws.Range("AB11").ClearContents() | Out-Null
$ws.Range("AF11").ClearContents() | Out-Null

# Row 12
$ws.Range("A12").Value = 112044164
$ws.Range("B12").Value = 88924
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 256703
$ws.Range("F12").Value = "Tallfingersvamp"
$ws.Range("G12").Value = "Ramaria eosanguinea"
$ws.Range("H12").Value = "R.H.Petersen"
$ws.Range("Q12").Value = 554725
$ws.Range("R12").Value = 6697591
$ws.Range("Z12").ClearContents() | Out-Null
$ws.Range("AB12").ClearContents() | Out-Null
# AF12: target is an empty string cell; leaving blank (engine cannot materialize empty cells)

# Row 23
$ws.Range("A23").Value = 112044157
$ws.Range("B23").Value = 89405
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 1202
$ws.Range("F23").Value = "Ullticka"
$ws.Range("G23").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H23").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q23").Value = 554764
$ws.Range("R23").Value = 6697616
$ws.Range("Z23").ClearContents() | Out-Null
$ws.Range("AB23").ClearContents() | Out-Null

# Row 24
$ws.Range("A24").Value = 112044169
$ws.Range("B24").Value = 89845
$ws.Range("D24").Value = "VU"
$ws.Range("E24").Value = 1209
$ws.Range("F24").Value = "Rynkskinn"
$ws.Range("G24").Value = "Phlebia centrifuga"
$ws.Range("H24").Value = "P.Karst."
$ws.Range("Q24").Value = 554765
$ws.Range("R24").Value = 6697617
$ws.Range("Z24").ClearContents() | Out-Null
$ws.Range("AB24").ClearContents() | Out-Null

# Row 25
$ws.Range("A25").Value = 112044179
$ws.Range("B25").Value = 96348
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("Q25").Value = 554795
$ws.Range("R25").Value = 6697596
$ws.Range("Z25").ClearContents() | Out-Null
$ws.Range("AB25").ClearContents() | Out-Null

# Row 26
$ws.Range("A26").Value = 112044168
$ws.Range("B26").Value = 89845
$ws.Range("E26").Value = 1209
$ws.Range("F26").Value = "Rynkskinn"
$ws.Range("G26").Value = "Phlebia centrifuga"
$ws.Range("H26").Value = "P.Karst."
$ws.Range("Q26").Value = 554761
$ws.Range("R26").Value = 6697614
$ws.Range("Z26").ClearContents() | Out-Null
$ws.Range("AB26").ClearContents() | Out-Null

# Row 27
$ws.Range("Q27").Value = 554839
$ws.Range("R27").Value = 6697581
$ws.Range("Z27").ClearContents() | Out-Null
$ws.Range("AB27").ClearContents() | Out-Null

# Row 28
$ws.Range("A28").Value = 112044155
$ws.Range("Q28").Value = 554761
$ws.Range("R28").Value = 6697629
$ws.Range("Z28").ClearContents() | Out-Null
$ws.Range("AB28").ClearContents() | Out-Null

# Row 29
$ws.Range("A29").Value = 112044156
$ws.Range("B29").Value = 89405
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 1202
$ws.Range("F29").Value = "Ullticka"
$ws.Range("G29").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H29").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q29").Value = 554762
$ws.Range("R29").Value = 6697614
$ws.Range("Z29").ClearContents() | Out-Null
$ws.Range("AB29").ClearContents() | Out-Null

# Row 30
$ws.Range("A30").Value = 112044154
$ws.Range("Q30").Value = 554768
$ws.Range("R30").Value = 6697637
$ws.Range("Z30").ClearContents() | Out-Null
$ws.Range("AB30").ClearContents() | Out-Null
